$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newF = @{
    2 = -1; 3 = 2; 4 = 0; 5 = 4; 6 = 4; 7 = -1; 8 = 3; 9 = 2; 10 = -2;
    11 = 6; 12 = -4; 13 = 0; 14 = 6; 15 = 3; 16 = 3; 17 = 2; 18 = -2;
    19 = -6; 20 = 4; 21 = -2; 22 = -5; 23 = 2; 24 = -1; 25 = 1; 26 = -4;
    27 = -5; 28 = -1; 29 = 3; 30 = 6; 31 = -5; 32 = 0; 33 = 3; 34 = 2;
    35 = 0; 36 = -2; 37 = -2; 38 = 4; 39 = -1; 40 = 0
}

foreach ($row in $newF.Keys) {
    $ws.Range("F$row").Value = $newF[$row]
}
